# Workbook "Hortaliza, Vega Central Mapocho de Santiago - Perejil"
# A new daily price record was inserted into the data table at row 534
# (pushing the existing rows 534:604 down to 535:605), growing the sheet
# from A1:R604 to A1:R605.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 534; rows 534:604 shift down to 535:605
$ws.Rows.Item(534).Insert()

# Fill in the values for the newly inserted row 534
$ws.Cells.Item(534, 1).Value  = 9
$ws.Cells.Item(534, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(534, 3).Value  = "Metropolitana"
$ws.Cells.Item(534, 4).Value  = 45154
$ws.Cells.Item(534, 5).Value  = 13
$ws.Cells.Item(534, 6).Value  = 100112044
$ws.Cells.Item(534, 7).Value  = "Perejil"
$ws.Cells.Item(534, 8).Value  = "Sin especificar"
$ws.Cells.Item(534, 9).Value  = "Primera"
$ws.Cells.Item(534, 10).Value = 70
$ws.Cells.Item(534, 11).Value = 12000
$ws.Cells.Item(534, 12).Value = 14000
$ws.Cells.Item(534, 13).Value = 13000
$ws.Cells.Item(534, 14).Value = "`$/docena de atados"
$ws.Cells.Item(534, 15).Value = "Región Metropolitana"
$ws.Cells.Item(534, 16).Value = 4333
$ws.Cells.Item(534, 17).Value = 3
$ws.Cells.Item(534, 18).Value = "Hortaliza"
